$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.8985665440559387
$ws.Range("B1").Value = 1.504457116127014
$ws.Range("C1").Value = 4.167180061340332
$ws.Range("D1").Value = 1.307035803794861
$ws.Range("E1").Value = 0.8377922773361206
